$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5032.32142456955
$ws.Range("C2").Value = 5115.18150386138
$ws.Range("F2").Value = 66.4730551788262

$ws.Range("B3").Value = 5054.91188051526
$ws.Range("C3").Value = 5507.48458985367
$ws.Range("F3").Value = 168.969793805767

$ws.Range("B4").Value = 5109.59008243961
$ws.Range("C4").Value = 4891.21995626175
$ws.Range("F4").Value = 144.308954950922

$ws.Range("B5").Value = 1360.12056462293
$ws.Range("C5").Value = 2977.76475728495
$ws.Range("F5").Value = 50.5848459859174

$ws.Range("B6").Value = 1120.5712356947
$ws.Range("C6").Value = 2806.55644652367
$ws.Range("F6").Value = 49.3323506178739

$ws.Range("B7").Value = 5232.97164748688
$ws.Range("C7").Value = 5416.81585191815
$ws.Range("F7").Value = 177.610039476303

$ws.Range("C9").Value = 5115.72188555132
$ws.Range("F9").Value = 135.897816605606

$ws.Range("C10").Value = 5043.29021650307
$ws.Range("F10").Value = 132.879830395262

$ws.Range("C11").Value = 4408.00610298797
$ws.Range("F11").Value = 106.409658998799

$ws.Range("C12").Value = 2475.58003566172
$ws.Range("F12").Value = 11.6648342008471

$ws.Range("C13").Value = 2219.53815276753
$ws.Range("F13").Value = 0.645145905706651

$ws.Range("C14").Value = 4951.01776196653
$ws.Range("F14").Value = 132.614002583785

$ws.Range("C15").Value = 4970.55689616187
$ws.Range("F15").Value = 133.428133175258
